# Generate Report for Archive
#
# 1. Change the status text "Ready for handoff" -> "In Translation" everywhere
#    it appears (Overview!E2:F3, zh-cn!C2:C3, de-de!C2:C3).
# 2. Narrow the (now shorter) status columns from ~17.22 characters to
#    ~13.41 characters: Overview columns E:F, and column C on the zh-cn /
#    de-de sheets. (12.42 is the nominal "ColumnWidth" input that this
#    host's column-width quantizer resolves closest to the 13.41-character
#    target stored width.)

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

if ($wsOverview.Range("E2").Value2 -eq $oldStatus) { $wsOverview.Range("E2").Value = $newStatus }
if ($wsOverview.Range("F2").Value2 -eq $oldStatus) { $wsOverview.Range("F2").Value = $newStatus }
if ($wsOverview.Range("E3").Value2 -eq $oldStatus) { $wsOverview.Range("E3").Value = $newStatus }
if ($wsOverview.Range("F3").Value2 -eq $oldStatus) { $wsOverview.Range("F3").Value = $newStatus }

$wsOverview.Range("E1:F1").EntireColumn.ColumnWidth = 12.42

# --- zh-cn sheet --------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

if ($wsZhCn.Range("C2").Value2 -eq $oldStatus) { $wsZhCn.Range("C2").Value = $newStatus }
if ($wsZhCn.Range("C3").Value2 -eq $oldStatus) { $wsZhCn.Range("C3").Value = $newStatus }

$wsZhCn.Range("C1").EntireColumn.ColumnWidth = 12.42

# --- de-de sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

if ($wsDeDe.Range("C2").Value2 -eq $oldStatus) { $wsDeDe.Range("C2").Value = $newStatus }
if ($wsDeDe.Range("C3").Value2 -eq $oldStatus) { $wsDeDe.Range("C3").Value = $newStatus }

$wsDeDe.Range("C1").EntireColumn.ColumnWidth = 12.42
